$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.834565043449402
$ws.Range("B1").Value = 2.508942604064941
$ws.Range("C1").Value = 4.768408298492432
$ws.Range("D1").Value = 4.604844570159912
$ws.Range("E1").Value = 1.262414455413818
